$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 17.76613666666667
$ws.Cells.Item(2, 8).Value = 53.29841
$ws.Cells.Item(2, 9).Value = 0.7872390387208499
$ws.Cells.Item(2, 10).Value = 0.7872390387208499
$ws.Cells.Item(2, 13).Value = 68.63737500000001
$ws.Cells.Item(2, 14).Value = 205.912125
$ws.Cells.Item(2, 15).Value = 0.5415701538216162
$ws.Cells.Item(2, 16).Value = 0.5415701538216162
$ws.Cells.Item(2, 17).Value = 1219.42098469125
$ws.Cells.Item(2, 18).Value = 10974.78886222125
$ws.Cells.Item(2, 19).Value = 0.4263451672944319
$ws.Cells.Item(2, 20).Value = 0.4263451672944319

$ws.Cells.Item(3, 7).Value = 17.76613666666667
$ws.Cells.Item(3, 8).Value = 53.29841
$ws.Cells.Item(3, 9).Value = 0.7872390387208499
$ws.Cells.Item(3, 10).Value = 0.7872390387208499
$ws.Cells.Item(3, 15).Value = 0.08718851262838957
$ws.Cells.Item(3, 16).Value = 0.08718851262838957
$ws.Cells.Item(3, 17).Value = 196.3171367048711
$ws.Cells.Item(3, 18).Value = 1766.85423034384
$ws.Cells.Item(3, 19).Value = 0.06863820086907409
$ws.Cells.Item(3, 20).Value = 0.06863820086907409

$ws.Cells.Item(4, 7).Value = 17.76613666666667
$ws.Cells.Item(4, 8).Value = 53.29841
$ws.Cells.Item(4, 9).Value = 0.7872390387208499
$ws.Cells.Item(4, 10).Value = 0.7872390387208499
$ws.Cells.Item(4, 13).Value = 16.21089566666667
$ws.Cells.Item(4, 14).Value = 48.632687
$ws.Cells.Item(4, 15).Value = 0.1279089892319285
$ws.Cells.Item(4, 16).Value = 0.1279089892319285
$ws.Cells.Item(4, 17).Value = 288.0049879030745
$ws.Cells.Item(4, 18).Value = 2592.04489112767
$ws.Cells.Item(4, 19).Value = 0.1006949497266989
$ws.Cells.Item(4, 20).Value = 0.1006949497266989

$ws.Cells.Item(5, 7).Value = 17.76613666666667
$ws.Cells.Item(5, 8).Value = 53.29841
$ws.Cells.Item(5, 9).Value = 0.7872390387208499
$ws.Cells.Item(5, 10).Value = 0.7872390387208499
$ws.Cells.Item(5, 13).Value = 20.32546233333333
$ws.Cells.Item(5, 14).Value = 60.976387
$ws.Cells.Item(5, 15).Value = 0.1603741949973873
$ws.Cells.Item(5, 16).Value = 0.1603741949973873
$ws.Cells.Item(5, 17).Value = 361.1049416271856
$ws.Cells.Item(5, 18).Value = 3249.94447464467
$ws.Cells.Item(5, 19).Value = 0.1262528271053733
$ws.Cells.Item(5, 20).Value = 0.1262528271053733

$ws.Cells.Item(6, 7).Value = 17.76613666666667
$ws.Cells.Item(6, 8).Value = 53.29841
$ws.Cells.Item(6, 9).Value = 0.7872390387208499
$ws.Cells.Item(6, 10).Value = 0.7872390387208499
$ws.Cells.Item(6, 13).Value = 10.513928
$ws.Cells.Item(6, 14).Value = 31.541784
$ws.Cells.Item(6, 15).Value = 0.08295814932067838
$ws.Cells.Item(6, 16).Value = 0.08295814932067838
$ws.Cells.Item(6, 17).Value = 186.7918817514933
$ws.Cells.Item(6, 18).Value = 1681.12693576344
$ws.Cells.Item(6, 19).Value = 0.06530789372527157
$ws.Cells.Item(6, 20).Value = 0.06530789372527157

$ws.Cells.Item(7, 9).Value = 0.03648413815195897
$ws.Cells.Item(7, 10).Value = 0.03648413815195897
$ws.Cells.Item(7, 13).Value = 68.63737500000001
$ws.Cells.Item(7, 14).Value = 205.912125
$ws.Cells.Item(7, 15).Value = 0.5415701538216162
$ws.Cells.Item(7, 16).Value = 0.5415701538216162
$ws.Cells.Item(7, 17).Value = 56.5133605965
$ws.Cells.Item(7, 18).Value = 508.6202453685
$ws.Cells.Item(7, 19).Value = 0.01975872031100552
$ws.Cells.Item(7, 20).Value = 0.01975872031100552

$ws.Cells.Item(8, 9).Value = 0.03648413815195897
$ws.Cells.Item(8, 10).Value = 0.03648413815195897
$ws.Cells.Item(8, 15).Value = 0.08718851262838957
$ws.Cells.Item(8, 16).Value = 0.08718851262838957
$ws.Cells.Item(8, 19).Value = 0.003180997739997985
$ws.Cells.Item(8, 20).Value = 0.003180997739997985

$ws.Cells.Item(9, 9).Value = 0.03648413815195897
$ws.Cells.Item(9, 10).Value = 0.03648413815195897
$ws.Cells.Item(9, 13).Value = 16.21089566666667
$ws.Cells.Item(9, 14).Value = 48.632687
$ws.Cells.Item(9, 15).Value = 0.1279089892319285
$ws.Cells.Item(9, 16).Value = 0.1279089892319285
$ws.Cells.Item(9, 17).Value = 13.34742467063422
$ws.Cells.Item(9, 18).Value = 120.126822035708
$ws.Cells.Item(9, 19).Value = 0.004666649234015112
$ws.Cells.Item(9, 20).Value = 0.004666649234015112

$ws.Cells.Item(10, 9).Value = 0.03648413815195897
$ws.Cells.Item(10, 10).Value = 0.03648413815195897
$ws.Cells.Item(10, 13).Value = 20.32546233333333
$ws.Cells.Item(10, 14).Value = 60.976387
$ws.Cells.Item(10, 15).Value = 0.1603741949973873
$ws.Cells.Item(10, 16).Value = 0.1603741949973873
$ws.Cells.Item(10, 17).Value = 16.73519976738978
$ws.Cells.Item(10, 18).Value = 150.616797906508
$ws.Cells.Item(10, 19).Value = 0.005851114286293887
$ws.Cells.Item(10, 20).Value = 0.005851114286293887

$ws.Cells.Item(11, 9).Value = 0.03648413815195897
$ws.Cells.Item(11, 10).Value = 0.03648413815195897
$ws.Cells.Item(11, 13).Value = 10.513928
$ws.Cells.Item(11, 14).Value = 31.541784
$ws.Cells.Item(11, 15).Value = 0.08295814932067838
$ws.Cells.Item(11, 16).Value = 0.08295814932067838
$ws.Cells.Item(11, 17).Value = 8.656761776650667
$ws.Cells.Item(11, 18).Value = 77.910855989856
$ws.Cells.Item(11, 19).Value = 0.003026656580646471
$ws.Cells.Item(11, 20).Value = 0.003026656580646471

$ws.Cells.Item(12, 7).Value = 3.885299333333334
$ws.Cells.Item(12, 8).Value = 11.655898
$ws.Cells.Item(12, 9).Value = 0.1721623203571791
$ws.Cells.Item(12, 10).Value = 0.172162320357179
$ws.Cells.Item(12, 13).Value = 68.63737500000001
$ws.Cells.Item(12, 14).Value = 205.912125
$ws.Cells.Item(12, 15).Value = 0.5415701538216162
$ws.Cells.Item(12, 16).Value = 0.5415701538216162
$ws.Cells.Item(12, 17).Value = 266.67674732925
$ws.Cells.Item(12, 18).Value = 2400.09072596325
$ws.Cells.Item(12, 19).Value = 0.09323797431812383
$ws.Cells.Item(12, 20).Value = 0.09323797431812382

$ws.Cells.Item(13, 7).Value = 3.885299333333334
$ws.Cells.Item(13, 8).Value = 11.655898
$ws.Cells.Item(13, 9).Value = 0.1721623203571791
$ws.Cells.Item(13, 10).Value = 0.172162320357179
$ws.Cells.Item(13, 15).Value = 0.08718851262838957
$ws.Cells.Item(13, 16).Value = 0.08718851262838957
$ws.Cells.Item(13, 17).Value = 42.93284773568356
$ws.Cells.Item(13, 18).Value = 386.395629621152
$ws.Cells.Item(13, 19).Value = 0.01501057664259476
$ws.Cells.Item(13, 20).Value = 0.01501057664259476

$ws.Cells.Item(14, 7).Value = 3.885299333333334
$ws.Cells.Item(14, 8).Value = 11.655898
$ws.Cells.Item(14, 9).Value = 0.1721623203571791
$ws.Cells.Item(14, 10).Value = 0.172162320357179
$ws.Cells.Item(14, 13).Value = 16.21089566666667
$ws.Cells.Item(14, 14).Value = 48.632687
$ws.Cells.Item(14, 15).Value = 0.1279089892319285
$ws.Cells.Item(14, 16).Value = 0.1279089892319285
$ws.Cells.Item(14, 17).Value = 62.98418212643624
$ws.Cells.Item(14, 18).Value = 566.8576391379261
$ws.Cells.Item(14, 19).Value = 0.02202110838071024
$ws.Cells.Item(14, 20).Value = 0.02202110838071024

$ws.Cells.Item(15, 7).Value = 3.885299333333334
$ws.Cells.Item(15, 8).Value = 11.655898
$ws.Cells.Item(15, 9).Value = 0.1721623203571791
$ws.Cells.Item(15, 10).Value = 0.172162320357179
$ws.Cells.Item(15, 13).Value = 20.32546233333333
$ws.Cells.Item(15, 14).Value = 60.976387
$ws.Cells.Item(15, 15).Value = 0.1603741949973873
$ws.Cells.Item(15, 16).Value = 0.1603741949973873
$ws.Cells.Item(15, 17).Value = 78.97050525339179
$ws.Cells.Item(15, 18).Value = 710.734547280526
$ws.Cells.Item(15, 19).Value = 0.0276103935361649
$ws.Cells.Item(15, 20).Value = 0.0276103935361649

$ws.Cells.Item(16, 7).Value = 3.885299333333334
$ws.Cells.Item(16, 8).Value = 11.655898
$ws.Cells.Item(16, 9).Value = 0.1721623203571791
$ws.Cells.Item(16, 10).Value = 0.172162320357179
$ws.Cells.Item(16, 13).Value = 10.513928
$ws.Cells.Item(16, 14).Value = 31.541784
$ws.Cells.Item(16, 15).Value = 0.08295814932067838
$ws.Cells.Item(16, 16).Value = 0.08295814932067838
$ws.Cells.Item(16, 17).Value = 40.84975744911467
$ws.Cells.Item(16, 18).Value = 367.647817042032
$ws.Cells.Item(16, 19).Value = 0.01428226747958533
$ws.Cells.Item(16, 20).Value = 0.01428226747958533

$ws.Cells.Item(17, 7).Value = 0.09285466666666665
$ws.Cells.Item(17, 8).Value = 0.278564
$ws.Cells.Item(17, 9).Value = 0.004114502770011991
$ws.Cells.Item(17, 10).Value = 0.004114502770011991
$ws.Cells.Item(17, 13).Value = 68.63737500000001
$ws.Cells.Item(17, 14).Value = 205.912125
$ws.Cells.Item(17, 15).Value = 0.5415701538216162
$ws.Cells.Item(17, 16).Value = 0.5415701538216162
$ws.Cells.Item(17, 17).Value = 6.3733005765
$ws.Cells.Item(17, 18).Value = 57.35970518849999
$ws.Cells.Item(17, 19).Value = 0.00222829189805486
$ws.Cells.Item(17, 20).Value = 0.00222829189805486

$ws.Cells.Item(18, 7).Value = 0.09285466666666665
$ws.Cells.Item(18, 8).Value = 0.278564
$ws.Cells.Item(18, 9).Value = 0.004114502770011991
$ws.Cells.Item(18, 10).Value = 0.004114502770011991
$ws.Cells.Item(18, 15).Value = 0.08718851262838957
$ws.Cells.Item(18, 16).Value = 0.08718851262838957
$ws.Cells.Item(18, 17).Value = 1.026050999815111
$ws.Cells.Item(18, 18).Value = 9.234458998335999
$ws.Cells.Item(18, 19).Value = 0.0003587373767227343
$ws.Cells.Item(18, 20).Value = 0.0003587373767227343

$ws.Cells.Item(19, 7).Value = 0.09285466666666665
$ws.Cells.Item(19, 8).Value = 0.278564
$ws.Cells.Item(19, 9).Value = 0.004114502770011991
$ws.Cells.Item(19, 10).Value = 0.004114502770011991
$ws.Cells.Item(19, 13).Value = 16.21089566666667
$ws.Cells.Item(19, 14).Value = 48.632687
$ws.Cells.Item(19, 15).Value = 0.1279089892319285
$ws.Cells.Item(19, 16).Value = 0.1279089892319285
$ws.Cells.Item(19, 17).Value = 1.505257313496444
$ws.Cells.Item(19, 18).Value = 13.547315821468
$ws.Cells.Item(19, 19).Value = 0.0005262818905042038
$ws.Cells.Item(19, 20).Value = 0.0005262818905042038

$ws.Cells.Item(20, 7).Value = 0.09285466666666665
$ws.Cells.Item(20, 8).Value = 0.278564
$ws.Cells.Item(20, 9).Value = 0.004114502770011991
$ws.Cells.Item(20, 10).Value = 0.004114502770011991
$ws.Cells.Item(20, 13).Value = 20.32546233333333
$ws.Cells.Item(20, 14).Value = 60.976387
$ws.Cells.Item(20, 15).Value = 0.1603741949973873
$ws.Cells.Item(20, 16).Value = 0.1603741949973873
$ws.Cells.Item(20, 17).Value = 1.887314029807555
$ws.Cells.Item(20, 18).Value = 16.985826268268
$ws.Cells.Item(20, 19).Value = 0.0006598600695551934
$ws.Cells.Item(20, 20).Value = 0.0006598600695551934

$ws.Cells.Item(21, 7).Value = 0.09285466666666665
$ws.Cells.Item(21, 8).Value = 0.278564
$ws.Cells.Item(21, 9).Value = 0.004114502770011991
$ws.Cells.Item(21, 10).Value = 0.004114502770011991
$ws.Cells.Item(21, 13).Value = 10.513928
$ws.Cells.Item(21, 14).Value = 31.541784
$ws.Cells.Item(21, 15).Value = 0.08295814932067838
$ws.Cells.Item(21, 16).Value = 0.08295814932067838
$ws.Cells.Item(21, 17).Value = 0.9762672797973332
$ws.Cells.Item(21, 18).Value = 8.786405518176
$ws.Cells.Item(21, 19).Value = 0.0003413315351749995
$ws.Cells.Item(21, 20).Value = 0.0003413315351749995
